$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that receive new values in this update (columns D and E).
# We must preserve them as literal text (matching the workbook's existing
# inline-string "General"-formatted cells) rather than let Excel
# auto-convert numeric-looking / percent-looking text into real numbers.
$updates = @{
    'D2' = '288.45'
    'E2' = '-4.46%'
    'D3' = '30.41'
    'E3' = '-6.25%'
    'D4' = '4.955'
    'E4' = '-2.18%'
    'D5' = '0.07247'
    'E5' = '-5.94%'
    'D6' = '1.794'
    'E6' = '-15.64%'
    'D7' = '7.571'
    'E7' = '-3.75%'
    'D8' = '3.711'
    'E8' = '-1.93%'
    'D9' = '0.9030'
    'E9' = '-2.15%'
    'D10' = '0.1665'
    'E10' = '-5.66%'
    'D11' = '0.07775'
    'E11' = '-2.05%'
    'D12' = '0.07924'
    'E12' = '-6.47%'
    'D13' = '0.03051'
    'E13' = '-0.03%'
    'E14' = '0.18%'
    'D16' = '0.005812'
    'E16' = '0.26%'
    'D18' = '3.470'
    'E18' = '0.20%'
    'D19' = '2.068'
    'E19' = '-3.91%'
    'D20' = '0.3303'
    'E20' = '-1.35%'
    'D21' = '0.1295'
    'E21' = '-1.68%'
    'D22' = '3.974'
    'E22' = '-6.83%'
    'D23' = '0.2098'
    'E23' = '6.57%'
    'D24' = '0.04496'
    'E24' = '-0.53%'
    'D25' = '0.001214'
    'E25' = '-1.41%'
    'E26' = '-4.59%'
    'D27' = '0.0001300'
    'E27' = '4.33%'
    'D39' = '0.01566'
    'E39' = '-9.13%'
    'D40' = '0.04332'
    'E40' = '-7.76%'
    'D41' = '0.007319'
    'E41' = '-3.26%'
    'D43' = '0.1305'
    'E43' = '-4.10%'
    'D44' = '0.002011'
    'E44' = '-9.69%'
    'D45' = '0.009429'
    'E45' = '-10.80%'
    'D46' = '0.00005877'
    'E46' = '-5.61%'
    'E47' = '0.30%'
    'E48' = '166.56%'
    'E49' = '0.30%'
    'D50' = '0.00002100'
    'E50' = '0.30%'
    'D51' = '0.0002000'
    'E51' = '0.30%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    # Force text storage so values like "288.45" or "-4.46%" stay as strings
    $range.NumberFormat = "@"
    $range.Value = $updates[$cellRef]
    # Remove the temporary text format again so the cell keeps the
    # workbook default (General) formatting, same as before the edit.
    $range.ClearFormats()
}
